# Insert a new row "E_STD / Stunden / Hours" immediately above the last
# row (row 28, "E_TEUR / 1 000 EUR / 1.000 EUR"), pushing that last row
# down to row 29 - matching the xl/worksheets/sheet1.xml diff:
#   dimension A1:C28 -> A1:C29
#   new row 28: E_STD | Stunden | Hours
#   old row 28 (E_TEUR | 1 000 EUR | 1.000 EUR) becomes row 29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 28
$newRow = $lastRow + 1

# Push the existing last row (and its formatting/style) down one row,
# using Copy so the cell style (s="4") is preserved faithfully.
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy($ws.Range("A" + $newRow + ":C" + $newRow))

# Overwrite the now-vacated row with the new unit entry.
$ws.Range("A" + $lastRow).Value = "E_STD"
$ws.Range("B" + $lastRow).Value = "Stunden"
$ws.Range("C" + $lastRow).Value = "Hours"

Write-Output "Inserted E_STD row at row $lastRow; moved E_TEUR row to row $newRow."
